# Auto-generated edit script: apply market-data value updates
# to the Atomos_Profits-derived Leve workbook (8 sheets: ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, WVR). Only literal numeric cell values
# change - no formulas are present in this workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2260.0144
$ws.Range("I15").Value = 2260.0144
$ws.Range("K15").Value = 6780.0432
$ws.Range("M15").Value = -6611.0432

$ws.Range("H74").Value = 3862.3572
$ws.Range("I74").Value = 3851.5
$ws.Range("J74").Value = 3864.1667
$ws.Range("K74").Value = 3851.5
$ws.Range("L74").Value = 3864.1667
$ws.Range("M74").Value = -2915.5
$ws.Range("N74").Value = -5736.1667

$ws.Range("H77").Value = 3862.3572
$ws.Range("I77").Value = 3851.5
$ws.Range("J77").Value = 3864.1667
$ws.Range("K77").Value = 19257.5
$ws.Range("L77").Value = 19320.8335
$ws.Range("M77").Value = -14577.5
$ws.Range("N77").Value = -28680.8335

$ws.Range("H106").Value = 2128.5293
$ws.Range("J106").Value = 2709.4
$ws.Range("L106").Value = 2709.4
$ws.Range("N106").Value = -3971.4

$ws.Range("H138").Value = 3586.5625
$ws.Range("I138").Value = 1346.9025
$ws.Range("J138").Value = 7579
$ws.Range("K138").Value = 4040.7075
$ws.Range("L138").Value = 22737
$ws.Range("M138").Value = 1099.2925
$ws.Range("N138").Value = -33017

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3776.45
$ws.Range("I32").Value = 3413.6814
$ws.Range("K32").Value = 3413.6814
$ws.Range("M32").Value = -3126.6814

$ws.Range("H38").Value = 1764.25
$ws.Range("I38").Value = 1764.25
$ws.Range("K38").Value = 1764.25
$ws.Range("M38").Value = -1297.25

$ws.Range("H39").Value = 3625
$ws.Range("I39").Value = 3625
$ws.Range("K39").Value = 3625
$ws.Range("M39").Value = -3105

$ws.Range("H40").Value = 42015.5
$ws.Range("I40").Value = 4000
$ws.Range("J40").Value = 80031
$ws.Range("K40").Value = 4000
$ws.Range("L40").Value = 80031
$ws.Range("M40").Value = -3824
$ws.Range("N40").Value = -80383

$ws.Range("H74").Value = 1446.2222
$ws.Range("I74").Value = 1060.3
$ws.Range("K74").Value = 1060.3
$ws.Range("M74").Value = -186.3

$ws.Range("H77").Value = 1446.2222
$ws.Range("I77").Value = 1060.3
$ws.Range("K77").Value = 5301.5
$ws.Range("M77").Value = -933.5

$ws.Range("H122").Value = 1545.5
$ws.Range("I122").Value = 1361.963
$ws.Range("J122").Value = 1875.8667
$ws.Range("K122").Value = 4085.889
$ws.Range("L122").Value = 5627.6001
$ws.Range("M122").Value = -1635.889
$ws.Range("N122").Value = -10527.6001

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

$ws.Range("H132").Value = 2737.5454
$ws.Range("I132").Value = 2696.4146
$ws.Range("J132").Value = 2858
$ws.Range("K132").Value = 8089.2438
$ws.Range("L132").Value = 8574
$ws.Range("M132").Value = -5559.2438
$ws.Range("N132").Value = -13634

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 33000
$ws.Range("J122").Value = 33000
$ws.Range("L122").Value = 33000
$ws.Range("N122").Value = -42800

$ws.Range("H137").Value = 30392.334
$ws.Range("J137").Value = 30272.637
$ws.Range("L137").Value = 30272.637
$ws.Range("N137").Value = -40472.637

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 19333.166
$ws.Range("J63").Value = 27875
$ws.Range("L63").Value = 83625
$ws.Range("N63").Value = -85123

$ws.Range("H66").Value = 19333.166
$ws.Range("J66").Value = 27875
$ws.Range("L66").Value = 250875
$ws.Range("N66").Value = -258363

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6834.4165
$ws.Range("I122").Value = 6201
$ws.Range("J122").Value = 7286.857
$ws.Range("K122").Value = 18603
$ws.Range("L122").Value = 21860.571
$ws.Range("M122").Value = -16153
$ws.Range("N122").Value = -26760.571

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = $null

$ws.Range("H132").Value = 3234.261
$ws.Range("I132").Value = 3206.3928
$ws.Range("J132").Value = 3277.611
$ws.Range("K132").Value = 9619.178400000001
$ws.Range("L132").Value = 9832.832999999999
$ws.Range("M132").Value = -7089.178400000001
$ws.Range("N132").Value = -14892.833

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1732.1818
$ws.Range("I40").Value = 1631.75
$ws.Range("K40").Value = 1631.75
$ws.Range("M40").Value = -1495.75

$ws.Range("H61").Value = 2438.7
$ws.Range("J61").Value = 6833.3335
$ws.Range("L61").Value = 6833.3335
$ws.Range("N61").Value = -7237.3335

$ws.Range("H113").Value = 2438.7
$ws.Range("J113").Value = 6833.3335
$ws.Range("L113").Value = 6833.3335
$ws.Range("N113").Value = -11173.3335

$ws.Range("H122").Value = 2628.375
$ws.Range("I122").Value = 2576.32
$ws.Range("J122").Value = 2814.2856
$ws.Range("K122").Value = 7728.960000000001
$ws.Range("L122").Value = 8442.856800000001
$ws.Range("M122").Value = -5278.960000000001
$ws.Range("N122").Value = -13342.8568

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5234.8335
$ws.Range("I62").Value = 6000
$ws.Range("J62").Value = 4852.25
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 4852.25
$ws.Range("M62").Value = -5376
$ws.Range("N62").Value = -6100.25

$ws.Range("H65").Value = 5234.8335
$ws.Range("I65").Value = 6000
$ws.Range("J65").Value = 4852.25
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 24261.25
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -30501.25

$ws.Range("H69").Value = 19635.5
$ws.Range("J69").Value = 19635.5
$ws.Range("L69").Value = 19635.5
$ws.Range("N69").Value = -21133.5

$ws.Range("H72").Value = 19635.5
$ws.Range("J72").Value = 19635.5
$ws.Range("L72").Value = 58906.5
$ws.Range("N72").Value = -66394.5

$ws.Range("H122").Value = 3750
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 4984.85
$ws.Range("I132").Value = 2236.2
$ws.Range("J132").Value = 13230.8
$ws.Range("K132").Value = 6708.599999999999
$ws.Range("L132").Value = 39692.39999999999
$ws.Range("M132").Value = -4178.599999999999
$ws.Range("N132").Value = -44752.39999999999
